$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right above the current row 380, pushing the
# existing rows 380..482 down to 382..484 (carries formatting forward,
# matching the canonical diff: dimension grows from R482 to R484).
$ws.Rows("380:381").Insert()

# Populate the two newly-inserted rows with the new weekly records.
$ws.Cells.Item(380, 1).Value = 6
$ws.Cells.Item(380, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(380, 3).Value = "Metropolitana"
$ws.Cells.Item(380, 4).Value = 44551
$ws.Cells.Item(380, 5).Value = 13
$ws.Cells.Item(380, 6).Value = 100112017
$ws.Cells.Item(380, 7).Value = "Apio"
$ws.Cells.Item(380, 8).Value = "Americana (o)"
$ws.Cells.Item(380, 9).Value = "Primera"
$ws.Cells.Item(380, 10).Value = 1900
$ws.Cells.Item(380, 11).Value = 6000
$ws.Cells.Item(380, 12).Value = 7000
$ws.Cells.Item(380, 13).Value = 6474
$ws.Cells.Item(380, 14).Value = "$/docena de matas"
$ws.Cells.Item(380, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(380, 16).Value = 1079
$ws.Cells.Item(380, 17).Value = 6
$ws.Cells.Item(380, 18).Value = "Hortaliza"

$ws.Cells.Item(381, 1).Value = 6
$ws.Cells.Item(381, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(381, 3).Value = "Metropolitana"
$ws.Cells.Item(381, 4).Value = 44551
$ws.Cells.Item(381, 5).Value = 13
$ws.Cells.Item(381, 6).Value = 100112017
$ws.Cells.Item(381, 7).Value = "Apio"
$ws.Cells.Item(381, 8).Value = "Americana (o)"
$ws.Cells.Item(381, 9).Value = "Segunda"
$ws.Cells.Item(381, 10).Value = 800
$ws.Cells.Item(381, 11).Value = 5000
$ws.Cells.Item(381, 12).Value = 5000
$ws.Cells.Item(381, 13).Value = 5000
$ws.Cells.Item(381, 14).Value = "$/docena de matas"
$ws.Cells.Item(381, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(381, 16).Value = 833
$ws.Cells.Item(381, 17).Value = 6
$ws.Cells.Item(381, 18).Value = "Hortaliza"
